# Updated cryptos list values (Price / Volume(1h)) per upstream diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.360.75"
$ws.Range("E2").Value = "  +1.08%  "

$ws.Range("D3").Value = "1.668.70"
$ws.Range("E3").Value = "  +1.05%  "

$ws.Range("E4").Value = "  +0.96%  "

$ws.Range("D5").Value = "'219.74"
$ws.Range("E5").Value = "  +1.23%  "

$ws.Range("D6").Value = "'0.5353"
$ws.Range("E6").Value = "  +1.55%  "

$ws.Range("E7").Value = "  +0.94%  "

$ws.Range("D8").Value = "'0.2663"
$ws.Range("E8").Value = "  +2.51%  "

$ws.Range("D9").Value = "'0.06396"
$ws.Range("E9").Value = "  +1.36%  "

$ws.Range("D10").Value = "'20.93"
$ws.Range("E10").Value = "  +3.00%  "

$ws.Range("D11").Value = "'0.07858"
$ws.Range("E11").Value = "  +0.94%  "

$ws.Range("D12").Value = "'4.563"
$ws.Range("E12").Value = "  +1.02%  "

$ws.Range("D13").Value = "1.665.87"
$ws.Range("E13").Value = "  +0.56%  "

$ws.Range("D14").Value = "1.896.55"
$ws.Range("E14").Value = "  +0.99%  "

$ws.Range("D15").Value = "'0.5541"
$ws.Range("E15").Value = "  +1.35%  "

$ws.Range("D16").Value = "0.0₅8196"
$ws.Range("E16").Value = "  +0.16%  "

$ws.Range("D17").Value = "'66.00"
$ws.Range("E17").Value = "  +1.09%  "

$ws.Range("D18").Value = "26.386.33"
$ws.Range("E18").Value = "  +1.21%  "

$ws.Range("D20").Value = "'4.683"
$ws.Range("E20").Value = "  +2.39%  "

$ws.Range("D21").Value = "'195.22"
$ws.Range("E21").Value = "  +2.55%  "

$ws.Range("E22").Value = "  +2.04%  "

$ws.Range("D23").Value = "'6.041"
$ws.Range("E23").Value = "  +0.49%  "

$ws.Range("D25").Value = "'146.52"
$ws.Range("E25").Value = "  +2.13%  "

$ws.Range("E26").Value = "  -0.39%  "

$ws.Range("D27").Value = "'7.237"
$ws.Range("E27").Value = "  +0.30%  "

$ws.Range("E28").Value = "  +1.18%  "

$ws.Range("E29").Value = "  +4.12%  "

$ws.Range("D30").Value = "'0.05860"
$ws.Range("E30").Value = "  +0.91%  "

$ws.Range("D31").Value = "'1.286"
$ws.Range("E31").Value = "  +1.26%  "

$ws.Range("D32").Value = "'3.584"
$ws.Range("E32").Value = "  +1.24%  "

$ws.Range("E33").Value = "  +1.14%  "

$ws.Range("D34").Value = "'1.615"
$ws.Range("E34").Value = "  +1.62%  "

$ws.Range("E35").Value = "  +3.30%  "

$ws.Range("D36").Value = "'2.838"
$ws.Range("E36").Value = "  +1.62%  "

$ws.Range("D37").Value = "'2.424"
$ws.Range("E37").Value = "  +0.49%  "

$ws.Range("D38").Value = "'0.5832"

$ws.Range("E39").Value = "  +0.13%  "

$ws.Range("D40").Value = "1.077.08"
$ws.Range("E40").Value = "  +4.71%  "

$ws.Range("D41").Value = "'0.8642"
$ws.Range("E41").Value = "  +1.74%  "

$ws.Range("D42").Value = "'5.860"
$ws.Range("E42").Value = "  +2.58%  "

$ws.Range("E43").Value = "  +0.98%  "

$ws.Range("D44").Value = "'104.37"
$ws.Range("E44").Value = "  -0.05%  "

$ws.Range("D45").Value = "1.806.34"
$ws.Range("E45").Value = "  +0.75%  "

$ws.Range("D46").Value = "'58.06"
$ws.Range("E46").Value = "  +1.88%  "

$ws.Range("E48").Value = "  +1.49%  "

$ws.Range("D49").Value = "'8.057"
$ws.Range("E49").Value = "  +2.78%  "

$ws.Range("D50").Value = "0.0₈103"
$ws.Range("E50").Value = "  -7.32%  "

$ws.Range("E51").Value = "  +0.60%  "
